$wb = $excel.ActiveWorkbook

$wsAdd = $wb.Worksheets.Item("Add Devices Loop A")
$wsOther = $wb.Worksheets.Item("Other Devices Loop A")

# --- Add Devices Loop A: add new "DC Unit Loading Details Name" column (E1:E3) ---
$wsAdd.Range("A7").Copy()
$wsAdd.Range("E1").PasteSpecial(-4122) # xlPasteFormats
$wsAdd.Range("E1").Value = "DC Unit Loading Details Name"

$wsAdd.Range("A8").Copy()
$wsAdd.Range("E2").PasteSpecial(-4122)
$wsAdd.Range("E2").Value = "Current (DC Units)"

$wsAdd.Range("A8").Copy()
$wsAdd.Range("E3").PasteSpecial(-4122)
$wsAdd.Range("E3").Value = "Current (worst case)"

# --- Other Devices Loop A: same new column ---
$wsOther.Range("A7").Copy()
$wsOther.Range("E1").PasteSpecial(-4122)
$wsOther.Range("E1").Value = "DC Unit Loading Details Name"

$wsOther.Range("A8").Copy()
$wsOther.Range("E2").PasteSpecial(-4122)
$wsOther.Range("E2").Value = "Current (DC Units)"

$wsOther.Range("A8").Copy()
$wsOther.Range("E3").PasteSpecial(-4122)
$wsOther.Range("E3").Value = "Current (worst case)"

# --- Update sheet selections / active sheet ---
$wsAdd.Activate()
$wsAdd.Range("C8").Select()

$wsOther.Activate()
$wsOther.Range("E1:E3").Select()
